# Re-render with lisa edits to ch1 and 2
# Shift the "date" column (F) values forward by 8 days for rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 44497
$ws.Range("F3").Value = 44496
$ws.Range("F4").Value = 44495
$ws.Range("F5").Value = 44494
$ws.Range("F6").Value = 44493
$ws.Range("F7").Value = 44492
